$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16: this shifts the existing rows 16-70 down to
# 17-71 (carrying all their values/styles along), so rows 17-71 end up
# exactly matching the "shifted" data from the diff automatically.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with a new weekly entry. It has the
# same Mercado/Región/Categoría/Variedad/Calidad/Volumen/Precios/etc. as
# the (now shifted-down) row 17, only with a new, more recent date.
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44804
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 100112022
$ws.Range("G16").Value = "Arveja Verde"
$ws.Range("H16").Value = "Perfection"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 28000
$ws.Range("L16").Value = 30000
$ws.Range("M16").Value = 29000
$ws.Range("N16").Value = '$/malla 25 kilos'
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 1160
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"
